# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# Commit: Updated cryptos list on Sat Mar 23 15:30:55 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values look like plain numbers (e.g. "1.00", "0.140").
# Force those specific cells to Text format BEFORE assigning the value so Excel
# keeps the exact string (trailing zeros, precision) instead of silently
# re-interpreting the input as a floating point number.
$numericLookingCells = @(
    "D5",
    "D6",
    "D11",
    "D12",
    "D14",
    "D16",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values
$ws.Range("D2").Value = "64.827.77"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "3.388.52"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "560.11"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "175.24"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("D8").Value = "3.379.28"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +10.14%  "
$ws.Range("D11").Value = "0.632"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "54.51"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("D14").Value = "9.15"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "3.945.55"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "18.27"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "3.393.08"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "64.694.71"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "470.60"
$ws.Range("E22").Value = "  +16.33%  "
$ws.Range("D23").Value = "4.95"
$ws.Range("E23").Value = "  +15.71%  "
$ws.Range("D24").Value = "4.13"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "86.48"
$ws.Range("E25").Value = "  +4.79%  "
$ws.Range("D26").Value = "13.75"
$ws.Range("E26").Value = "  +5.45%  "
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +4.79%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "30.66"
$ws.Range("E30").Value = "  +5.16%  "
$ws.Range("D31").Value = "6.73"
$ws.Range("E31").Value = "  +3.78%  "
$ws.Range("D32").Value = "11.54"
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").Value = "579.71"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "60.00"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "0.140"
$ws.Range("E37").Value = "  -4.64%  "
$ws.Range("D38").Value = "35.93"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "0.0₃0755"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").Value = "3.099.72"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.134"
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.18"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "136.97"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "8.39"
$ws.Range("E51").Value = "  +4.25%  "
